$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap category labels for rows 13 and 14 (shared-string reorder: "Family..." now
# at the index previously used by "Discrimination...", and vice versa).
$ws.Range("A13").Value = "Family; children; childcare"
$ws.Range("A14").Value = "Discrimination; gender inequality; racism; LGBT"

# Re-run: Russia (col L) recomputed throughout; "All" aggregate (col B) recalculated
# accordingly. Rows 13/14 also swap all their other country values (category order swap).
$ws.Range("B2").Value = 0.261037795077051
$ws.Range("L2").Value = 0.263660744832401
$ws.Range("B3").Value = 0.19287692482063
$ws.Range("L3").Value = 0.20355238386636
$ws.Range("B4").Value = 0.188436644013223
$ws.Range("L4").Value = 0.084485900521279
$ws.Range("B5").Value = 0.165026496197311
$ws.Range("L5").Value = 0.130835356619311
$ws.Range("B6").Value = 0.13881451950498
$ws.Range("L6").Value = 0.136912364546507
$ws.Range("B7").Value = 0.138273185837646
$ws.Range("L7").Value = 0.280111468081168
$ws.Range("B8").Value = 0.128826472505465
$ws.Range("L8").Value = 0.0459741097958687
$ws.Range("B9").Value = 0.127106283277647
$ws.Range("F9").Value = 0.15677241305404
$ws.Range("L9").Value = 0.134193145267982
$ws.Range("B10").Value = 0.114455094191174
$ws.Range("L10").Value = 0.0625074943637807
$ws.Range("B11").Value = 0.0900921929978382
$ws.Range("L11").Value = 0.0393939362920809
$ws.Range("B12").Value = 0.0827364469588425
$ws.Range("E12").Value = 0.0781231459534224
$ws.Range("L12").Value = 0.0690848019997287
$ws.Range("B13").Value = 0.0782534856915513
$ws.Range("C13").Value = 0.0702572220672497
$ws.Range("D13").Value = 0.0598568390567794
$ws.Range("E13").Value = 0.0498847795322759
$ws.Range("F13").Value = 0.0725313941240774
$ws.Range("G13").Value = 0.0713429670037117
$ws.Range("H13").Value = 0.0595647178068743
$ws.Range("I13").Value = 0.113866367038356
$ws.Range("J13").Value = 0.0556416961094614
$ws.Range("K13").Value = 0.0861912889898547
$ws.Range("L13").Value = 0.0688046498164585
$ws.Range("M13").Value = 0.11442732319241
$ws.Range("N13").Value = 0.0850570554794627
$ws.Range("B14").Value = 0.0782025760992126
$ws.Range("C14").Value = 0.0801606130053347
$ws.Range("D14").Value = 0.0853143000988915
$ws.Range("E14").Value = 0.0854716178134182
$ws.Range("F14").Value = 0.0845599979045629
$ws.Range("G14").Value = 0.0325015769502991
$ws.Range("H14").Value = 0.0680430990380986
$ws.Range("I14").Value = 0.0981141721056102
$ws.Range("J14").Value = 0.100751012095868
$ws.Range("K14").Value = 0.0596500618249719
$ws.Range("L14").Value = 0.0264130970851678
$ws.Range("M14").Value = 0.100448138308641
$ws.Range("N14").Value = 0.103404378607846
$ws.Range("B15").Value = 0.0719975965782097
$ws.Range("L15").Value = 0.0266693031101935
$ws.Range("B16").Value = 0.0678510222497633
$ws.Range("L16").Value = 0.0374756892607604
$ws.Range("B17").Value = 0.0671132616697566
$ws.Range("L17").Value = 0.0362911051776244
$ws.Range("B18").Value = 0.057396255145867
$ws.Range("L18").Value = 0.0545332898958408
$ws.Range("B19").Value = 0.0519052111375239
$ws.Range("L19").Value = 0.0699679981949483
$ws.Range("M19").Value = 0.0426027285771167
$ws.Range("B20").Value = 0.042448802314565
$ws.Range("C20").Value = 0.0623696912729769
$ws.Range("L20").Value = 0.00427594286072039
$ws.Range("B21").Value = 0.0400453044786933
$ws.Range("K21").Value = 0.0166432031700155
$ws.Range("L21").Value = 0.0137349231870967
$ws.Range("B22").Value = 0.0350372005822907
$ws.Range("C22").Value = 0.0349985902063247
$ws.Range("D22").Value = 0.0238260416955817
$ws.Range("L22").Value = 0.0385146174271064
$ws.Range("B23").Value = 0.0317780452883306
$ws.Range("L23").Value = 0.00131220261070944
$ws.Range("B24").Value = 0.0309368248291205
$ws.Range("L24").Value = 0.0312969682333077
$ws.Range("B25").Value = 0.0262552506904307
$ws.Range("L25").Value = 0.00453310360539263
$ws.Range("B26").Value = 0.0191186000044473
$ws.Range("L26").Value = 0.0302216684254651
$ws.Range("B27").Value = 0.0112777986161369
$ws.Range("F27").Value = 0.0070095978599204
$ws.Range("L27").Value = 0.00376730220275702
$ws.Range("B28").Value = 0.00847546908801203
$ws.Range("L28").Value = 0.0018486238275502
